$d = $word.ActiveDocument

# The end of the document currently has, right after the
# "LOT2053: Microbiologia (Requisito fraco)" paragraph:
#   1. a blank separator paragraph
#   2. "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3. "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github
#       pages. Original theme under Creative Commons Attribution"
#   4. another blank paragraph
#   5. a (blank) page-break paragraph
#
# The footer block (1)-(3) above must be removed, leaving just the blank
# paragraph (4) and the page-break paragraph (5) in place.

$jupiterText = "Ver no Jupiter Salvar em pdf Salvar em docx"
$copyrightText = "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"

# Locate the two footer paragraphs by searching for their full text, so the
# match range starts exactly at the beginning of each paragraph.
$rJupiter = $d.Content
$foundJupiter = $rJupiter.Find.Execute($jupiterText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$rCopyright = $d.Content
$foundCopyright = $rCopyright.Find.Execute($copyrightText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($foundJupiter -and $foundCopyright) {
    # Translate the found text ranges into paragraph indices.
    $jupiterIndex = 0
    $copyrightIndex = 0
    $i = 0
    foreach ($p in $d.Paragraphs) {
        $i = $i + 1
        $pStart = $p.Range.Start
        if ($pStart -eq $rJupiter.Start) { $jupiterIndex = $i }
        if ($pStart -eq $rCopyright.Start) { $copyrightIndex = $i }
    }

    if ($jupiterIndex -gt 1 -and $copyrightIndex -ge $jupiterIndex) {
        # Paragraph right before "Ver no Jupiter ..." is the blank
        # separator paragraph that must also be removed.
        $blankBefore = $d.Paragraphs.Item($jupiterIndex - 1)
        $copyrightPara = $d.Paragraphs.Item($copyrightIndex)

        $deleteStart = $blankBefore.Range.Start
        $deleteEnd = $copyrightPara.Range.End

        $deleteRange = $d.Range($deleteStart, $deleteEnd)
        $deleteRange.Delete()
    }
}
